$d = $word.ActiveDocument

# Locate the unique paragraph containing "...maslees<comment>c_019r_01</comment> &amp;<lb/>"
# by anchoring on "aslees", which occurs exactly once in the document, then
# expanding the found range out to its enclosing paragraph.
$anchor = $d.Content.Duplicate
$anchor.Find.Execute("aslees") | Out-Null
$anchor.Expand(4) | Out-Null   # wdParagraph

$paraStart = $anchor.Start
$paraEnd = $anchor.End

# First edit: the run rendered as " &" (Arial, 22pt) loses its trailing "&",
# becoming just " ".
$r1 = $d.Range($paraStart, $paraEnd)
$r1.Find.Execute(" &", $true, $false, $false, $false, $false, $true, 1, $false, " ", 2) | Out-Null

# Second edit: the following run rendered as "amp;<lb/>" (Courier New, 18pt,
# gray) loses the leading "amp;" text, becoming just "<lb/>".
$r2 = $d.Range($paraStart, $paraEnd)
$r2.Find.Execute("amp;<lb/>", $true, $false, $false, $false, $false, $true, 1, $false, "<lb/>", 2) | Out-Null
